$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 366.1111
$ws.Range("I6").Value = 37
$ws.Range("K6").Value = 111
$ws.Range("M6").Value = 1

$ws.Range("H12").Value = 376.66666
$ws.Range("I12").Value = 30
$ws.Range("K12").Value = 30
$ws.Range("M12").Value = 140

$ws.Range("H17").Value = 2079.1765
$ws.Range("I17").Value = 900
$ws.Range("J17").Value = 2384.889
$ws.Range("K17").Value = 2700
$ws.Range("L17").Value = 7154.667
$ws.Range("M17").Value = -2532
$ws.Range("N17").Value = -7490.667

$ws.Range("H40").Value = 6884.6665
$ws.Range("I40").Value = 6188.6665
$ws.Range("J40").Value = 7232.6665
$ws.Range("K40").Value = 6188.6665
$ws.Range("L40").Value = 7232.6665
$ws.Range("M40").Value = -6013.6665
$ws.Range("N40").Value = -7582.6665

$ws.Range("H69").Value = 7541.8335
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 7541.8335
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 22625.5005
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -24373.5005

$ws.Range("H72").Value = 7541.8335
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 7541.8335
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 67876.5015
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -76612.5015

$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("N88").ClearContents()

$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("N91").ClearContents()

$ws.Range("H94").Value = 5501.273
$ws.Range("I94").Value = 4612.6665
$ws.Range("K94").Value = 4612.6665
$ws.Range("M94").Value = -4161.6665

$ws.Range("H99").Value = 2310.5
$ws.Range("I99").Value = 331
$ws.Range("J99").Value = 4290
$ws.Range("K99").Value = 993
$ws.Range("L99").Value = 12870
$ws.Range("M99").Value = 505
$ws.Range("N99").Value = -15866

$ws.Range("H100").Value = 5118
$ws.Range("I100").Value = 4863.3335
$ws.Range("K100").Value = 4863.3335
$ws.Range("M100").Value = -4322.3335

$ws.Range("H135").Value = 2198.6
$ws.Range("I135").Value = 1915
$ws.Range("K135").Value = 17235
$ws.Range("M135").Value = -14700

$ws.Range("H137").Value = 2421.111
$ws.Range("I137").Value = 1339.8
$ws.Range("K137").Value = 4019.4
$ws.Range("M137").Value = -1469.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2210.5715
$ws.Range("I2").Value = 2163
$ws.Range("K2").Value = 2163
$ws.Range("M2").Value = -2050

$ws.Range("H32").Value = 5000955.5
$ws.Range("I32").Value = 1005.7895
$ws.Range("K32").Value = 1005.7895
$ws.Range("M32").Value = -718.7895

$ws.Range("H45").Value = 2746.7693
$ws.Range("I45").Value = 2086.4
$ws.Range("K45").Value = 2086.4
$ws.Range("M45").Value = -1709.4

$ws.Range("H97").Value = 947.8570999999999
$ws.Range("I97").Value = 944.5454999999999
$ws.Range("J97").Value = 960
$ws.Range("K97").Value = 944.5454999999999
$ws.Range("L97").Value = 960
$ws.Range("M97").Value = -448.5454999999999
$ws.Range("N97").Value = -1952

$ws.Range("H116").Value = 2210.5715
$ws.Range("I116").Value = 2163
$ws.Range("K116").Value = 2163
$ws.Range("M116").Value = 131

$ws.Range("H132").Value = 1365.625
$ws.Range("I132").Value = 1313.0952
$ws.Range("K132").Value = 3939.2856
$ws.Range("M132").Value = -1409.2856

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2210.5715
$ws.Range("I3").Value = 2163
$ws.Range("K3").Value = 2163
$ws.Range("M3").Value = -2049

$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 200
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 200
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -27
$ws.Range("N22").ClearContents()

$ws.Range("H70").Value = 250000
$ws.Range("J70").Value = 250000
$ws.Range("L70").Value = 250000
$ws.Range("N70").Value = -250586

$ws.Range("H73").Value = 250000
$ws.Range("J73").Value = 250000
$ws.Range("L73").Value = 250000
$ws.Range("N73").Value = -252028

$ws.Range("H86").Value = 6791.8184
$ws.Range("I86").Value = 3740.6
$ws.Range("K86").Value = 3740.6
$ws.Range("M86").Value = -2617.6

$ws.Range("H89").Value = 6791.8184
$ws.Range("I89").Value = 3740.6
$ws.Range("K89").Value = 18703
$ws.Range("M89").Value = -13087

$ws.Range("H134").Value = 4291
$ws.Range("I134").Value = 1115.6666
$ws.Range("J134").Value = 16198.5
$ws.Range("K134").Value = 3346.9998
$ws.Range("L134").Value = 48595.5
$ws.Range("M134").Value = -811.9998000000001
$ws.Range("N134").Value = -53665.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4958
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()

$ws.Range("H132").Value = 1284.2222
$ws.Range("I132").Value = 1284.2222
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3852.6666
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1322.6666
$ws.Range("N132").ClearContents()

$ws.Range("H136").Value = 4958
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 806.25
$ws.Range("J34").Value = 894.5714
$ws.Range("L34").Value = 2683.7142
$ws.Range("N34").Value = -2851.7142

$ws.Range("H57").Value = 2337.5
$ws.Range("I57").Value = 1675
$ws.Range("J57").Value = 3000
$ws.Range("K57").Value = 5025
$ws.Range("L57").Value = 9000
$ws.Range("M57").Value = -4466
$ws.Range("N57").Value = -10118

$ws.Range("H86").Value = 610.1111
$ws.Range("I86").Value = 611.5
$ws.Range("J86").Value = 599
$ws.Range("K86").Value = 1834.5
$ws.Range("L86").Value = 1797
$ws.Range("M86").Value = -648.5
$ws.Range("N86").Value = -4169

$ws.Range("H89").Value = 610.1111
$ws.Range("I89").Value = 611.5
$ws.Range("J89").Value = 599
$ws.Range("K89").Value = 5503.5
$ws.Range("L89").Value = 5391
$ws.Range("M89").Value = 424.5
$ws.Range("N89").Value = -17247

$ws.Range("H92").Value = 2880.625
$ws.Range("I92").Value = 998.6667
$ws.Range("J92").Value = 4009.8
$ws.Range("K92").Value = 2996.0001
$ws.Range("L92").Value = 12029.4
$ws.Range("M92").Value = -1748.0001
$ws.Range("N92").Value = -14525.4

$ws.Range("H126").Value = 1166.6666
$ws.Range("I126").Value = 1166.6666
$ws.Range("K126").Value = 3499.9998
$ws.Range("M126").Value = 1440.0002

$ws.Range("H137").Value = 5774.5
$ws.Range("J137").Value = 5774.5
$ws.Range("L137").Value = 17323.5
$ws.Range("N137").Value = -27523.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 1395
$ws.Range("I70").Value = 1395
$ws.Range("K70").Value = 1395
$ws.Range("M70").Value = -1125

$ws.Range("H73").Value = 1395
$ws.Range("I73").Value = 1395
$ws.Range("K73").Value = 1395
$ws.Range("M73").Value = -459

$ws.Range("H122").Value = 2928.818
$ws.Range("I122").Value = 2707.5
$ws.Range("J122").Value = 3519
$ws.Range("K122").Value = 8122.5
$ws.Range("L122").Value = 10557
$ws.Range("M122").Value = -5672.5
$ws.Range("N122").Value = -15457

$ws.Range("H132").Value = 36334.45
$ws.Range("I132").Value = 38803.668
$ws.Range("K132").Value = 116411.004
$ws.Range("M132").Value = -113881.004
